# Rename the "T2" condition labels back to their base names in the
# header row of Sheet1 (A1:D1). The "space" label (used elsewhere in
# column D) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "square"
$ws.Range("B1").Value = "loc1"
$ws.Range("C1").Value = "loc2"
$ws.Range("D1").Value = "corrAns"
